# Insert a new "finished" column before column AE (old AE -> AF),
# fill it with 1 for every data row, and set header style/text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at AE, shifting the old AE (recoding_done) to AF.
$ws.Range("AE1").EntireColumn.Insert()

# Header text/style for the new column (match the other header cells' style,
# i.e. same look as W1:AB1 which use the Arial 10 font style).
$ws.Range("W1").Copy()
$ws.Range("AE1").PasteSpecial(-4122)
$ws.Range("AE1").Value = "finished"

# Fill data rows 2-101 with 1 in the new column (plain, unstyled cells).
$ws.Range("AE2:AE101").Style = "Standard"
for ($r = 2; $r -le 101; $r++) {
    $ws.Cells.Item($r, 31).Value = 1
}

# Update view state to match target (scrolled/selected one column to the right).
$ws.Application.ActiveWindow.ScrollColumn = 25
$ws.Range("AD8").Select()
